$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 770, pushing existing rows 770-878 down to 771-879
$ws.Rows.Item(770).Insert()

# Populate the newly inserted row 770 with the new record
$ws.Range("A770").Value = 10
$ws.Range("B770").Value = 'Vega Modelo de Temuco'
$ws.Range("C770").Value = 'La Araucanía'
$ws.Range("D770").Value = 45131
$ws.Range("E770").Value = 9
$ws.Range("F770").Value = 100112032
$ws.Range("G770").Value = 'Zapallo italiano'
$ws.Range("H770").Value = 'Sin especificar'
$ws.Range("I770").Value = 'Primera'
$ws.Range("J770").Value = 295
$ws.Range("K770").Value = 17000
$ws.Range("L770").Value = 18000
$ws.Range("M770").Value = 17627
$ws.Range("N770").Value = '$/caja 50 unidades'
$ws.Range("O770").Value = 'Región de Arica y Parinacota'
$ws.Range("P770").Value = 353
$ws.Range("Q770").Value = 50
$ws.Range("R770").Value = 'Hortaliza'
